$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.575558
$ws.Range("H2").Value = 10.726674
$ws.Range("I2").Value = 0.025194653521236
$ws.Range("J2").Value = 0.02519465352123599
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 4.901647840660666
$ws.Range("R2").Value = 44.114830565946
$ws.Range("S2").Value = 0.0002779129789079086
$ws.Range("T2").Value = 0.0002779129789079085

# Row 3
$ws.Range("G3").Value = 3.575558
$ws.Range("H3").Value = 10.726674
$ws.Range("I3").Value = 0.025194653521236
$ws.Range("J3").Value = 0.02519465352123599
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 332.2184103498446
$ws.Range("R3").Value = 2989.965693148602
$ws.Range("S3").Value = 0.01883607535051537
$ws.Range("T3").Value = 0.01883607535051537

# Row 4
$ws.Range("G4").Value = 3.575558
$ws.Range("H4").Value = 10.726674
$ws.Range("I4").Value = 0.025194653521236
$ws.Range("J4").Value = 0.02519465352123599
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 106.260163214072
$ws.Range("R4").Value = 956.3414689266478
$ws.Range("S4").Value = 0.006024724635069454
$ws.Range("T4").Value = 0.006024724635069451

# Row 5
$ws.Range("G5").Value = 3.575558
$ws.Range("H5").Value = 10.726674
$ws.Range("I5").Value = 0.025194653521236
$ws.Range("J5").Value = 0.02519465352123599
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 0.9866430500779999
$ws.Range("R5").Value = 8.879787450701999
$ws.Range("S5").Value = 0.0000559405567432612
$ws.Range("T5").Value = 0.00005594055674326118

# Row 6
$ws.Range("I6").Value = 0.7460690747908298
$ws.Range("J6").Value = 0.7460690747908298
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 145.1485675859681
$ws.Range("R6").Value = 1306.337108273713
$ws.Range("S6").Value = 0.008229614226344597
$ws.Range("T6").Value = 0.008229614226344596

# Row 7
$ws.Range("I7").Value = 0.7460690747908298
$ws.Range("J7").Value = 0.7460690747908298
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.5577775974416317
$ws.Range("T7").Value = 0.5577775974416316

# Row 8
$ws.Range("I8").Value = 0.7460690747908298
$ws.Range("J8").Value = 0.7460690747908298
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 3146.597018666049
$ws.Range("R8").Value = 28319.37316799444
$ws.Range("S8").Value = 0.1784053402666233
$ws.Range("T8").Value = 0.1784053402666233

# Row 9
$ws.Range("I9").Value = 0.7460690747908298
$ws.Range("J9").Value = 0.7460690747908298
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 29.21666959619233
$ws.Range("R9").Value = 262.950026365731
$ws.Range("S9").Value = 0.001656522856230227
$ws.Range("T9").Value = 0.001656522856230227

# Row 10
$ws.Range("G10").Value = 32.36130266666667
$ws.Range("H10").Value = 97.08390800000001
$ws.Range("I10").Value = 0.2280292497513723
$ws.Range("J10").Value = 0.2280292497513723
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 44.36334394157022
$ws.Range("R10").Value = 399.2700954741321
$ws.Range("S10").Value = 0.002515306988568995
$ws.Range("T10").Value = 0.002515306988568995

# Row 11
$ws.Range("G11").Value = 32.36130266666667
$ws.Range("H11").Value = 97.08390800000001
$ws.Range("I11").Value = 0.2280292497513723
$ws.Range("J11").Value = 0.2280292497513723
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 3006.809155038231
$ws.Range("R11").Value = 27061.28239534409
$ws.Range("S11").Value = 0.1704796665220275
$ws.Range("T11").Value = 0.1704796665220274

# Row 12
$ws.Range("G12").Value = 32.36130266666667
$ws.Range("H12").Value = 97.08390800000001
$ws.Range("I12").Value = 0.2280292497513723
$ws.Range("J12").Value = 0.2280292497513723
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 961.7288555184907
$ws.Range("R12").Value = 8655.559699666417
$ws.Range("S12").Value = 0.05452797504579858
$ws.Range("T12").Value = 0.05452797504579857

# Row 13
$ws.Range("G13").Value = 32.36130266666667
$ws.Range("H13").Value = 97.08390800000001
$ws.Range("I13").Value = 0.2280292497513723
$ws.Range("J13").Value = 0.2280292497513723
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 8.929810219142666
$ws.Range("R13").Value = 80.368291972284
$ws.Range("S13").Value = 0.0005063011949772642
$ws.Range("T13").Value = 0.0005063011949772641

# Row 14
$ws.Range("G14").Value = 0.1003386666666667
$ws.Range("H14").Value = 0.301016
$ws.Range("I14").Value = 0.000707021936561918
$ws.Range("J14").Value = 0.0007070219365619179
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 0.1375519034515556
$ws.Range("R14").Value = 1.237967131064
$ws.Range("S14").Value = 0.000007798899571194482
$ws.Range("T14").Value = 0.000007798899571194481

# Row 15
$ws.Range("G15").Value = 0.1003386666666667
$ws.Range("H15").Value = 0.301016
$ws.Range("I15").Value = 0.000707021936561918
$ws.Range("J15").Value = 0.0007070219365619179
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 9.32283921464089
$ws.Range("R15").Value = 83.905552931768
$ws.Range("S15").Value = 0.0005285851008160345
$ws.Range("T15").Value = 0.0005285851008160344

# Row 16
$ws.Range("G16").Value = 0.1003386666666667
$ws.Range("H16").Value = 0.301016
$ws.Range("I16").Value = 0.000707021936561918
$ws.Range("J16").Value = 0.0007070219365619179
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 2.981913059914667
$ws.Range("R16").Value = 26.837217539232
$ws.Range("S16").Value = 0.0001690681110239825
$ws.Range("T16").Value = 0.0001690681110239825

# Row 17
$ws.Range("G17").Value = 0.1003386666666667
$ws.Range("H17").Value = 0.301016
$ws.Range("I17").Value = 0.000707021936561918
$ws.Range("J17").Value = 0.0007070219365619179
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 8.929810219142666
$ws.Range("R17").Value = 80.368291972284
$ws.Range("S17").Value = 0.0005063011949772642
$ws.Range("T17").Value = 0.0005063011949772641
